# Daily attendance processing - 2025-10-29 05:22:41
#
# The "Recorded By" column (G) lists the users who recorded / touched each
# attendance session, separated by ", ". For every row whose G cell holds
# more than one name, flip the order of the comma-separated names (the most
# recently-recorded name moves to the front). Rows with only a single name
# (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val.Contains(",")) {
        $parts = $val.Split(",")
        $count = $parts.Count

        $reversed = ""
        for ($i = $count - 1; $i -ge 0; $i--) {
            $piece = $parts[$i].Trim()
            if ($i -eq ($count - 1)) {
                $reversed = $piece
            } else {
                $reversed = $reversed + ", " + $piece
            }
        }

        $cell.Value = $reversed
    }
}
